$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are stored as text in the source sheet (prices use "." as a
# thousands separator and percentages keep padding spaces), so any cell
# whose new value would otherwise be auto-detected as a number is forced
# back to Text format before the value is written.

$ws.Range("D2").Value = "36.699.78"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.130.15"
$ws.Range("E3").Value = "  +11.12%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.59"
$ws.Range("E5").Value = "  +2.87%  "
$ws.Range("E6").Value = "  -3.98%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.05"
$ws.Range("E8").Value = "  +6.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.04"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "2.441.05"
$ws.Range("E13").Value = "  +11.11%  "
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.842"
$ws.Range("E15").Value = "  +5.64%  "
$ws.Range("D16").Value = "2.129.52"
$ws.Range("E16").Value = "  +11.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.13"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "36.743.17"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.61"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.35"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.58"
$ws.Range("E22").Value = "  -3.81%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  -7.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.92"
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.82"
$ws.Range("E27").Value = "  +16.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.25"
$ws.Range("E28").Value = "  +5.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.03"
$ws.Range("E29").Value = "  -7.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.33"
$ws.Range("E30").Value = "  +58.89%  "
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0962"
$ws.Range("E33").Value = "  +13.69%  "
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("E35").Value = "  +17.58%  "
$ws.Range("E36").Value = "  +9.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.17"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  -9.20%  "
$ws.Range("E41").Value = "  +8.99%  "
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.05"
$ws.Range("E43").Value = "  -7.11%  "
$ws.Range("E44").Value = "  +11.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.14"
$ws.Range("E45").Value = "  -5.89%  "
$ws.Range("D46").Value = "1.358.63"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.27"
$ws.Range("E47").Value = "  +13.13%  "
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("D49").Value = "2.326.39"
$ws.Range("E49").Value = "  +11.73%  "
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("E51").Value = "  -3.10%  "
